# Actualización automática 2025-09-10 15:50:09
#
# Updates a few raw sales figures for client "MADECOR-HOME CENTER S.A.S."
# (PORCELANATO) and "ROMERO RODAS SILVIA MARELIS" (INODOROS / LAVABOS),
# then propagates the change through the dependent summary sheets
# (monthly totals, compliance totals and "x de 53" counters).

$wb  = $excel.ActiveWorkbook
$wsGrupo  = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl   = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO: raw per-client / per-category sale values
# ---------------------------------------------------------------------
$wsGrupo.Range("M27").Value = 165.89          # MADECOR-HOME CENTER S.A.S. / PORCELANATO
$wsGrupo.Range("H46").Value = 1128.6          # ROMERO RODAS SILVIA MARELIS / INODOROS
$wsGrupo.Range("I46").Value = 172.8           # ROMERO RODAS SILVIA MARELIS / LAVABOS

# Row 55 holds "<count> de 53" counters (number of clients with sales > 0
# in that category). The three edits above turn a previously-zero cell
# into a positive one, so each affected counter increases by one.
$wsGrupo.Range("H55").Value = "2 de 53"
$wsGrupo.Range("I55").Value = "6 de 53"
$wsGrupo.Range("M55").Value = "8 de 53"

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL: "septiembre" column (F) is the per-client sum across
#    all categories in VENTAS POR GRUPO, and row 59 is the column total.
# ---------------------------------------------------------------------
$wsMensual.Range("F27").Value = 623.8099999999999   # MADECOR-HOME CENTER S.A.S.
$wsMensual.Range("F46").Value = 1301.4              # ROMERO RODAS SILVIA MARELIS
$wsMensual.Range("F59").Value = 19447.16            # column total

# ---------------------------------------------------------------------
# 3) CUMPLIMIENTO MENSUAL: per-category VENTA (D), POR CUMPLIR (E = C-D)
#    and CUMPLIMIENTO (F = D/C), plus the TOTAL row 15.
# ---------------------------------------------------------------------
# Row 6 -> INODOROS
$wsCumpl.Range("D6").Value = 1926.9
$wsCumpl.Range("E6").Value = 980.6836814602598
$wsCumpl.Range("F6").Value = 0.6627152340572581

# Row 7 -> LAVABOS
$wsCumpl.Range("D7").Value = 1050.3
$wsCumpl.Range("E7").Value = -163.5889837124259
$wsCumpl.Range("F7").Value = 1.184489625940738

# Row 12 -> PORCELANATO
$wsCumpl.Range("D12").Value = 9124.719999999999
$wsCumpl.Range("E12").Value = 52739.0003947566
$wsCumpl.Range("F12").Value = 0.1474971104514009

# Row 15 -> TOTAL
$wsCumpl.Range("D15").Value = 18926.75
$wsCumpl.Range("E15").Value = 103128.0855108344
$wsCumpl.Range("F15").Value = 0.1550675966321707

# Column E widened slightly to fit the new (longer) numbers.
# (Excel's ColumnWidth COM property is offset ~0.83 from the raw OOXML
# <col width> attribute for this workbook's default font; 23.17 here
# serialises to width="24" in the saved file, matching the target.)
$wsCumpl.Range("E1").EntireColumn.ColumnWidth = 23.17
